$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: -25C data
$ws.Range("B6").Value = 10165
$ws.Range("C6").Value = 4.0232200000000002
$ws.Range("D6").Value = 3.1842600000000001
$ws.Range("E6").Value = -0.00105
$ws.Range("F6").Value = -4.9984299999999999

# Row 7: -30C data
$ws.Range("B7").Value = 15333
$ws.Range("C7").Value = 3.9992000000000001
$ws.Range("D7").Value = 3.0104799999999998
$ws.Range("E7").Value = -0.00104
$ws.Range("F7").Value = -4.9983899999999997

# Update selection to match diff (active cell F8)
$ws.Range("F8").Select()

$wb.Save()
